$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column O (Extracted Objects), shifting
# old O..U (Extracted Objects .. Result String) to Q..W.
$ws.Range("O1:P1").EntireColumn.Insert()

# Rename M1/N1 headers
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# New O1/P1 headers
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Fill new O/P columns (rows 2-6) with copies of the M/N values for that row
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item($r, 13).Value2
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($r, 14).Value2
}
